$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.047.31"
$ws.Range("E2").Value = "  +1.55%  "
$ws.Range("D3").Value = "1.637.67"
$ws.Range("E3").Value = "  +2.36%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'215.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.60%  "
$ws.Range("E6").Value = "  +1.63%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").Value = "'29.52"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +9.77%  "
$ws.Range("E9").Value = "  +4.44%  "
$ws.Range("D10").Value = "'0.0615"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.44%  "
$ws.Range("D11").Value = "'0.0918"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.80%  "
$ws.Range("D12").Value = "1.871.29"
$ws.Range("E12").Value = "  +2.31%  "
$ws.Range("D13").Value = "1.647.39"
$ws.Range("E13").Value = "  +3.00%  "
$ws.Range("D14").Value = "'0.578"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.06%  "
$ws.Range("D15").Value = "'9.63"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +26.14%  "
$ws.Range("D16").Value = "'3.92"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.82%  "
$ws.Range("D17").Value = "30.034.47"
$ws.Range("E17").Value = "  +1.45%  "
$ws.Range("D18").Value = "'65.03"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.07%  "
$ws.Range("D19").Value = "'249.39"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.26%  "
$ws.Range("E20").Value = "  +2.30%  "
$ws.Range("D22").Value = "'4.21"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.49%  "
$ws.Range("D23").Value = "'9.90"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.02%  "
$ws.Range("D24").Value = "'2.13"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.77%  "
$ws.Range("D25").Value = "'159.51"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.90%  "
$ws.Range("D26").Value = "'15.76"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.65%  "
$ws.Range("D27").Value = "'0.112"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.04%  "
$ws.Range("D28").Value = "'6.69"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.47%  "
$ws.Range("E29").Value = "  -0.12%  "
$ws.Range("D30").Value = "'0.0493"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.13%  "
$ws.Range("D31").Value = "'1.14"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.88%  "
$ws.Range("D32").Value = "'3.43"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.42%  "
$ws.Range("D33").Value = "'3.23"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.17%  "
$ws.Range("D34").Value = "1.439.37"
$ws.Range("E34").Value = "  +1.26%  "
$ws.Range("E35").Value = "  +8.34%  "
$ws.Range("E36").Value = "  +1.75%  "
$ws.Range("E37").Value = "  +0.27%  "
$ws.Range("E38").Value = "  -0.10%  "
$ws.Range("E39").Value = "  +1.66%  "
$ws.Range("D40").Value = "'76.14"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +15.19%  "
$ws.Range("D41").Value = "'0.561"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.19%  "
$ws.Range("D42").Value = "'0.840"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.84%  "
$ws.Range("E43").Value = "  +1.42%  "
$ws.Range("E44").Value = "  +0.97%  "
$ws.Range("D45").Value = "'55.11"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.70%  "
$ws.Range("E46").Value = "  +4.91%  "
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("E48").Value = "  +2.17%  "
$ws.Range("D49").Value = "1.777.87"
$ws.Range("E49").Value = "  +2.35%  "
$ws.Range("D50").Value = "'90.57"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.09%  "
$ws.Range("E51").Value = "  +6.27%  "
